# Atualização de bases das ligas, do dia: 10-06-2024 às 21:53
# Swap the match-data (columns B:AD, i.e. everything except the row-index
# column A) between specific rows. Column A keeps its original value in
# place; only the "id" .. "PL_AhUnder" fields move between rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Read the current (pre-edit) row data for every row involved -----
$row107 = $ws.Range("B107:AD107").Value2
$row108 = $ws.Range("B108:AD108").Value2

$row128 = $ws.Range("B128:AD128").Value2
$row129 = $ws.Range("B129:AD129").Value2

$row143 = $ws.Range("B143:AD143").Value2
$row144 = $ws.Range("B144:AD144").Value2
$row145 = $ws.Range("B145:AD145").Value2

$row211 = $ws.Range("B211:AD211").Value2
$row212 = $ws.Range("B212:AD212").Value2

$row214 = $ws.Range("B214:AD214").Value2
$row215 = $ws.Range("B215:AD215").Value2

# --- Write back the swapped/rotated data ------------------------------

# Rows 107 <-> 108 : straight swap
$ws.Range("B107:AD107").Value2 = $row108
$ws.Range("B108:AD108").Value2 = $row107

# Rows 128 <-> 129 : straight swap
$ws.Range("B128:AD128").Value2 = $row129
$ws.Range("B129:AD129").Value2 = $row128

# Rows 143, 144, 145 : cyclic rotation (143<-145, 144<-143, 145<-144)
$ws.Range("B143:AD143").Value2 = $row145
$ws.Range("B144:AD144").Value2 = $row143
$ws.Range("B145:AD145").Value2 = $row144

# Rows 211 <-> 212 : straight swap
$ws.Range("B211:AD211").Value2 = $row212
$ws.Range("B212:AD212").Value2 = $row211

# Rows 214 <-> 215 : straight swap
$ws.Range("B214:AD214").Value2 = $row215
$ws.Range("B215:AD215").Value2 = $row214
